$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 70

# Text columns - use a leading apostrophe so Excel stores them as literal
# text instead of auto-converting to a date/number, then clear the
# resulting "quote prefix" formatting so the cell keeps the sheet's
# default (unstyled) look.
$ws.Cells.Item($row, 1).Value = "'2024-01-18"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "09:00:10"

$ws.Cells.Item($row, 3).Value = "Thursday"

$ws.Cells.Item($row, 4).Value = "'02"
$ws.Cells.Item($row, 4).ClearFormats()

# Numeric columns
$ws.Cells.Item($row, 5).Value = 138954
$ws.Cells.Item($row, 6).Value = 139690
$ws.Cells.Item($row, 7).Value = 170800
$ws.Cells.Item($row, 8).Value = 148623
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119686
$ws.Cells.Item($row, 11).Value = 222821
$ws.Cells.Item($row, 12).Value = 254582
$ws.Cells.Item($row, 13).Value = 185123
$ws.Cells.Item($row, 14).Value = 110347
$ws.Cells.Item($row, 15).Value = 41307
$ws.Cells.Item($row, 16).Value = 30930
$ws.Cells.Item($row, 17).Value = 73466
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42130
$ws.Cells.Item($row, 20).Value = -1
